$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted before the existing row 149
# (Fruta / Pomelo weekly data set), shifting all subsequent rows down
# by one and pushing the former last row (271) out to a new row 272.
$ws.Rows(149).Insert()

$ws.Range("A149").Value = 4
$ws.Range("B149").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C149").Value = 'Los Lagos'
$ws.Range("D149").Value = 44673
$ws.Range("E149").Value = 10
$ws.Range("F149").Value = 'Fruta'
$ws.Range("G149").Value = 100102
$ws.Range("H149").Value = 'Cítricos'
$ws.Range("I149").Value = 100102006
$ws.Range("J149").Value = 'Pomelo'
$ws.Range("K149").Value = 'Start Ruby'
$ws.Range("L149").Value = 'Primera'
$ws.Range("M149").Value = 120
$ws.Range("N149").Value = 14000
$ws.Range("O149").Value = 15000
$ws.Range("P149").Value = 14500
$ws.Range("Q149").Value = '$/caja 14 kilos empedrada'
$ws.Range("R149").Value = "Región de O'Higgins"
$ws.Range("S149").Value = 1036
$ws.Range("T149").Value = 14
